$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Overwrite data rows 2-5 (A2:AH5) with the new dataset values
$data = New-Object "object[,]" 4,34
$data[0,0] = 45097.50694444445
$data[0,1] = 7.205
$data[0,2] = 5.605
$data[0,3] = 1.192
$data[0,4] = 15.61
$data[0,5] = 12.374
$data[0,6] = 4.794
$data[0,7] = 14.798
$data[0,8] = 8.957
$data[0,9] = 4.429
$data[0,10] = 5.631
$data[0,11] = 6.249
$data[0,12] = 7.306
$data[0,13] = 2.788
$data[0,14] = 6.015
$data[0,15] = 7.738
$data[0,16] = 5.138
$data[0,17] = 0.492
$data[0,18] = 0.931
$data[0,19] = 84.14
$data[0,20] = 16.424
$data[0,21] = 5.552
$data[0,22] = 10.174
$data[0,23] = 6.283
$data[0,24] = 0.894
$data[0,25] = 9.741
$data[0,26] = 4.361
$data[0,27] = 4.885
$data[0,28] = 6.06
$data[0,29] = 8.134
$data[0,30] = 1.522
$data[0,31] = 13.18
$data[0,32] = 3.625
$data[0,33] = 6.387
$data[1,0] = 45097.51388888889
$data[1,1] = 5.005
$data[1,2] = 3.813
$data[1,3] = 0.567
$data[1,4] = 11.126
$data[1,5] = 8.677
$data[1,6] = 3.398
$data[1,7] = 15.294
$data[1,8] = 6.169
$data[1,9] = 3.216
$data[1,10] = 3.824
$data[1,11] = 4.409
$data[1,12] = 5.051
$data[1,13] = 1.855
$data[1,14] = 4.136
$data[1,15] = 5.456
$data[1,16] = 3.629
$data[1,17] = 0.191
$data[1,18] = 0.504
$data[1,19] = 55.623
$data[1,20] = 11.627
$data[1,21] = 3.817
$data[1,22] = 7.36
$data[1,23] = 4.345
$data[1,24] = 0.572
$data[1,25] = 8.85
$data[1,26] = 3.073
$data[1,27] = 3.363
$data[1,28] = 4.033
$data[1,29] = 5.333
$data[1,30] = 0.707
$data[1,31] = 14.501
$data[1,32] = 2.453
$data[1,33] = 4.471
$data[2,0] = 45097.52083333334
$data[2,1] = 24.305
$data[2,2] = 18.276
$data[2,3] = 1.066
$data[2,4] = 53.109
$data[2,5] = 43.488
$data[2,6] = 18.734
$data[2,7] = 69.486
$data[2,8] = 29.497
$data[2,9] = 13.604
$data[2,10] = 19.523
$data[2,11] = 21.268
$data[2,12] = 22.767
$data[2,13] = 6.534
$data[2,14] = 19.174
$data[2,15] = 27.068
$data[2,16] = 16.059
$data[2,17] = 0.281
$data[2,18] = 0.953
$data[2,19] = 284.528
$data[2,20] = 53.527
$data[2,21] = 17.698
$data[2,22] = 35.898
$data[2,23] = 19.222
$data[2,24] = 2.528
$data[2,25] = 35.593
$data[2,26] = 15.427
$data[2,27] = 14.021
$data[2,28] = 16.526
$data[2,29] = 22.788
$data[2,30] = 0.461
$data[2,31] = 63.413
$data[2,32] = 10.29
$data[2,33] = 21.919
$data[3,0] = 45097.52777777778
$data[3,1] = 21.95
$data[3,2] = 16.49
$data[3,3] = 0.92
$data[3,4] = 47.96
$data[3,5] = 39.3
$data[3,6] = 16.96
$data[3,7] = 68.06
$data[3,8] = 26.62
$data[3,9] = 12.27
$data[3,10] = 17.63
$data[3,11] = 19.21
$data[3,12] = 20.53
$data[3,13] = 5.85
$data[3,14] = 17.29
$data[3,15] = 24.48
$data[3,16] = 14.48
$data[3,17] = 0.23
$data[3,18] = 0.82
$data[3,19] = 255.93
$data[3,20] = 48.39
$data[3,21] = 15.96
$data[3,22] = 32.5
$data[3,23] = 17.33
$data[3,24] = 2.28
$data[3,25] = 33.7
$data[3,26] = 13.94
$data[3,27] = 12.62
$data[3,28] = 14.86
$data[3,29] = 20.5
$data[3,30] = 0.34
$data[3,31] = 62.14
$data[3,32] = 9.26
$data[3,33] = 19.8

$ws.Range("A2:AH5").Value = $data

# Row 6 no longer exists in the new dataset -> delete it (shifts dimension to AH5)
$ws.Rows.Item(6).Delete()

# Column-width tweaks from the custom-accuracy formatting pass
$ws.Columns.Item(2).ColumnWidth = 7.17
$ws.Columns.Item(3).ColumnWidth = 7.17
$ws.Columns.Item(7).ColumnWidth = 7.17
$ws.Columns.Item(10).ColumnWidth = 7.17
$ws.Columns.Item(11).ColumnWidth = 7.17
$ws.Columns.Item(12).ColumnWidth = 7.17
$ws.Columns.Item(13).ColumnWidth = 7.17
$ws.Columns.Item(15).ColumnWidth = 7.17
$ws.Columns.Item(16).ColumnWidth = 7.17
$ws.Columns.Item(17).ColumnWidth = 7.17
$ws.Columns.Item(20).ColumnWidth = 8.17
$ws.Columns.Item(22).ColumnWidth = 7.17
$ws.Columns.Item(24).ColumnWidth = 7.17
$ws.Columns.Item(26).ColumnWidth = 7.17
$ws.Columns.Item(27).ColumnWidth = 7.17
$ws.Columns.Item(28).ColumnWidth = 7.17
$ws.Columns.Item(29).ColumnWidth = 7.17
$ws.Columns.Item(30).ColumnWidth = 7.17
$ws.Columns.Item(34).ColumnWidth = 7.17
